$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-05-06 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-07 Tuesday", 2)

$d.Content.Find.Execute("15÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2)
$d.Content.Find.Execute("23÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷3=", 2)
$d.Content.Find.Execute("21÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷2=", 2)
$d.Content.Find.Execute("18÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷5=", 2)
$d.Content.Find.Execute("28÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "84÷3=", 2)
$d.Content.Find.Execute("26÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷3=", 2)
$d.Content.Find.Execute("89÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "64÷9=", 2)
$d.Content.Find.Execute("88÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2)
$d.Content.Find.Execute("20÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷6=", 2)
$d.Content.Find.Execute("24÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "29÷4=", 2)
$d.Content.Find.Execute("56÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷6=", 2)
$d.Content.Find.Execute("59÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=", 2)
$d.Content.Find.Execute("80÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷3=", 2)
$d.Content.Find.Execute("44÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "80÷6=", 2)
$d.Content.Find.Execute("28÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷9=", 2)
$d.Content.Find.Execute("26÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=", 2)
$d.Content.Find.Execute("98÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷7=", 2)
$d.Content.Find.Execute("62÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷4=", 2)
$d.Content.Find.Execute("18÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷5=", 2)
$d.Content.Find.Execute("14÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "28÷9=", 2)
$d.Content.Find.Execute("35÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷4=", 2)
$d.Content.Find.Execute("82÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "71÷6=", 2)
$d.Content.Find.Execute("56÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "39÷8=", 2)
$d.Content.Find.Execute("62÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "23÷9=", 2)
$d.Content.Find.Execute("57÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷5=", 2)
